$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 2557
    $ws.Range("F4").Value = 508
    $ws.Range("F6").Value = 6563
    $ws.Range("F7").Value = 392
    $ws.Range("F9").Value = 133
}
